{"js": "// Replace every paragraph's text (the title line plus all 100 arithmetic\n// table-cell entries) with the updated values from the commit, in document\n// order. Using Range.insertText(..., Word.InsertLocation.replace) rewrites\n// only the text of each run/paragraph while leaving paragraph/run\n// formatting (fonts, sizes, alignment, etc.) untouched.\nconst newTexts = [\n  \"2023-09-08 Friday\",\n  \"51-7=44\",\n  \"52+39=91\",\n  \"41-36=5\",\n  \"35+16=51\",\n  \"60-46=14\",\n  \"6+75=81\",\n  \"91-33=58\",\n  \"79+16=95\",\n  \"60-21=39\",\n  \"62-26=36\",\n  \"55+9=64\",\n  \"44-37=7\",\n  \"71-52=19\",\n  \"80-55=25\",\n  \"19+67=86\",\n  \"49+37=86\",\n  \"45+7=52\",\n  \"18+68=86\",\n  \"85-66=19\",\n  \"8+73=81\",\n  \"8+76=84\",\n  \"9+8=17\",\n  \"87-68=19\",\n  \"6+15=21\",\n  \"38+6=44\",\n  \"67+6=73\",\n  \"43-26=17\",\n  \"4+8=12\",\n  \"52-48=4\",\n  \"38+55=93\",\n  \"23+28=51\",\n  \"90-52=38\",\n  \"29+37=66\",\n  \"68+29=97\",\n  \"37+24=61\",\n  \"33-27=6\",\n  \"7+5=12\",\n  \"65+7=72\",\n  \"34+47=81\",\n  \"84-8=76\",\n  \"34+38=72\",\n  \"9+44=53\",\n  \"66+9=75\",\n  \"34+29=63\",\n  \"16+45=61\",\n  \"84-65=19\",\n  \"16+76=92\",\n  \"37+46=83\",\n  \"47+38=85\",\n  \"90-4=86\",\n  \"41-9=32\",\n  \"82-17=65\",\n  \"72-5=67\",\n  \"85-26=59\",\n  \"34-28=6\",\n  \"27+26=53\",\n  \"90-67=23\",\n  \"72-3=69\",\n  \"91-55=36\",\n  \"88+3=91\",\n  \"33+39=72\",\n  \"43-39=4\",\n  \"29+46=75\",\n  \"93-27=66\",\n  \"64+9=73\",\n  \"8+63=71\",\n  \"79+3=82\",\n  \"38+49=87\",\n  \"63-48=15\",\n  \"81-52=29\",\n  \"63-55=8\",\n  \"52-4=48\",\n  \"80-1=79\",\n  \"6+18=24\",\n  \"61-15=46\",\n  \"95-6=89\",\n  \"9+22=31\",\n  \"84-67=17\",\n  \"35+18=53\",\n  \"17+29=46\",\n  \"48+9=57\",\n  \"19+76=95\",\n  \"26+49=75\",\n  \"63-15=48\",\n  \"52+39=91\",\n  \"59+27=86\",\n  \"26+15=41\",\n  \"19+16=35\",\n  \"32-16=16\",\n  \"69+26=95\",\n  \"54+38=92\",\n  \"9+27=36\",\n  \"5+28=33\",\n  \"32+59=91\",\n  \"26+8=34\",\n  \"81-17=64\",\n  \"54-47=7\",\n  \"90-16=74\",\n  \"94-65=29\",\n  \"8+33=41\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newTexts.length) {\n  throw new Error(\n    `Expected ${newTexts.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the title/date line.\n$d = $word.ActiveDocument\n$d.Paragraphs.Item(1).Range.Text = '2023-09-08 Friday'\n\n# Update every cell of the single 20x5 arithmetic-answers table,\n# row by row, left to right, matching document order.\n$newValues = @(\n    '51-7=44',\n    '52+39=91',\n    '41-36=5',\n    '35+16=51',\n    '60-46=14',\n    '6+75=81',\n    '91-33=58',\n    '79+16=95',\n    '60-21=39',\n    '62-26=36',\n    '55+9=64',\n    '44-37=7',\n    '71-52=19',\n    '80-55=25',\n    '19+67=86',\n    '49+37=86',\n    '45+7=52',\n    '18+68=86',\n    '85-66=19',\n    '8+73=81',\n    '8+76=84',\n    '9+8=17',\n    '87-68=19',\n    '6+15=21',\n    '38+6=44',\n    '67+6=73',\n    '43-26=17',\n    '4+8=12',\n    '52-48=4',\n    '38+55=93',\n    '23+28=51',\n    '90-52=38',\n    '29+37=66',\n    '68+29=97',\n    '37+24=61',\n    '33-27=6',\n    '7+5=12',\n    '65+7=72',\n    '34+47=81',\n    '84-8=76',\n    '34+38=72',\n    '9+44=53',\n    '66+9=75',\n    '34+29=63',\n    '16+45=61',\n    '84-65=19',\n    '16+76=92',\n    '37+46=83',\n    '47+38=85',\n    '90-4=86',\n    '41-9=32',\n    '82-17=65',\n    '72-5=67',\n    '85-26=59',\n    '34-28=6',\n    '27+26=53',\n    '90-67=23',\n    '72-3=69',\n    '91-55=36',\n    '88+3=91',\n    '33+39=72',\n    '43-39=4',\n    '29+46=75',\n    '93-27=66',\n    '64+9=73',\n    '8+63=71',\n    '79+3=82',\n    '38+49=87',\n    '63-48=15',\n    '81-52=29',\n    '63-55=8',\n    '52-4=48',\n    '80-1=79',\n    '6+18=24',\n    '61-15=46',\n    '95-6=89',\n    '9+22=31',\n    '84-67=17',\n    '35+18=53',\n    '17+29=46',\n    '48+9=57',\n    '19+76=95',\n    '26+49=75',\n    '63-15=48',\n    '52+39=91',\n    '59+27=86',\n    '26+15=41',\n    '19+16=35',\n    '32-16=16',\n    '69+26=95',\n    '54+38=92',\n    '9+27=36',\n    '5+28=33',\n    '32+59=91',\n    '26+8=34',\n    '81-17=64',\n    '54-47=7',\n    '90-16=74',\n    '94-65=29',\n    '8+33=41'\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n"}
